# KANBAN proyecto CASAI.xlsx - update board state
#
# - Minimize the workbook window (bookViews/workbookView minimized flag).
# - Move the finished "Modificar generacion de presuúesto..." card from the
#   "EN PROCESO" column (B41) to the "TERMINADAS" column (C41).
# - Add a new "PENDIENTES" card in A42: Creacion de modulo de planeación "En pausa".
# - Add a new "EN PROCESO" card in B43 about the conciliacion export, and give
#   its still-empty "TERMINADAS" cell (C43) the green underlined Tahoma look
#   used for other in-flight/linked items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the Excel window (mirrors minimized="1" on workbookView).
$excel.WindowState = -4140

# --- Move the completed task from "EN PROCESO" to "TERMINADAS" -------------
$tareaTerminada = $ws.Range("B41").Text
$ws.Range("C41").Value = $tareaTerminada
$ws.Range("B41").Value = ""

# --- New pending task -------------------------------------------------------
$ws.Range("A42").Value = "Creacion de modulo de planeación ""En pausa"""

# --- New in-process task -----------------------------------------------------
$ws.Range("B43").Value = "Modificar archivo de conciliacion para que muestre el resumen de forma ordenada"

# Row 43 now wraps onto two lines just like the other long entries (row 41 etc.)
$ws.Rows.Item(43).RowHeight = 28.5

# Style the (still empty) C43 cell with the green underlined Tahoma font used
# elsewhere on the sheet for this kind of reference/placeholder cell.
$ws.Range("C43").Font.Name = "Tahoma"
$ws.Range("C43").Font.Size = 11
$ws.Range("C43").Font.Underline = $true
$ws.Range("C43").Font.Color = 5287936

# Keep the active selection in sync with the last-edited cell.
$ws.Range("B43").Select()
